# Refresh the "cryptos" price/volume snapshot (columns D = Price, E = Volume(1h))
# for the rows whose figures moved since the last scrape. Values that look
# numeric (e.g. "194.76") are written with a leading apostrophe so Excel
# keeps them as text, matching the original inline-string cells instead of
# silently re-typing them as numbers (which would mangle things like
# "0.0000301" -> 3.01E-05 or drop the trailing zero in "1.00").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '69.682.25'
$ws.Cells.Item(2, 5).Value = '  +0.19%  '
$ws.Cells.Item(3, 4).Value = '3.505.82'
$ws.Cells.Item(4, 4).Value = '''0.998'
$ws.Cells.Item(4, 5).Value = '  -0.13%  '
$ws.Cells.Item(5, 4).Value = '''605.13'
$ws.Cells.Item(5, 5).Value = '  -1.24%  '
$ws.Cells.Item(6, 4).Value = '''194.76'
$ws.Cells.Item(6, 5).Value = '  +2.72%  '
$ws.Cells.Item(7, 5).Value = '  +0.07%  '
$ws.Cells.Item(8, 5).Value = '  +0.02%  '
$ws.Cells.Item(9, 5).Value = '  -5.29%  '
$ws.Cells.Item(10, 5).Value = '  +0.39%  '
$ws.Cells.Item(11, 4).Value = '''53.45'
$ws.Cells.Item(11, 5).Value = '  +0.82%  '
$ws.Cells.Item(12, 4).Value = '''0.0000301'
$ws.Cells.Item(12, 5).Value = '  -2.03%  '
$ws.Cells.Item(13, 4).Value = '''9.52'
$ws.Cells.Item(13, 5).Value = '  +0.28%  '
$ws.Cells.Item(14, 4).Value = '4.064.55'
$ws.Cells.Item(14, 5).Value = '  +0.20%  '
$ws.Cells.Item(15, 4).Value = '''593.93'
$ws.Cells.Item(15, 5).Value = '  -0.82%  '
$ws.Cells.Item(16, 5).Value = '  +0.63%  '
$ws.Cells.Item(17, 4).Value = '69.839.43'
$ws.Cells.Item(17, 5).Value = '  +0.28%  '
$ws.Cells.Item(18, 4).Value = '''12.75'
$ws.Cells.Item(18, 5).Value = '  +1.35%  '
$ws.Cells.Item(19, 5).Value = '  +2.17%  '
$ws.Cells.Item(20, 4).Value = '3.505.82'
$ws.Cells.Item(20, 5).Value = '  +0.20%  '
$ws.Cells.Item(21, 4).Value = '''0.989'
$ws.Cells.Item(21, 5).Value = '  +0.27%  '
$ws.Cells.Item(22, 4).Value = '''18.30'
$ws.Cells.Item(22, 5).Value = '  +6.79%  '
$ws.Cells.Item(23, 4).Value = '''5.30'
$ws.Cells.Item(23, 5).Value = '  +3.20%  '
$ws.Cells.Item(24, 5).Value = '  -1.06%  '
$ws.Cells.Item(25, 4).Value = '''101.60'
$ws.Cells.Item(25, 5).Value = '  -3.82%  '
$ws.Cells.Item(26, 4).Value = '''3.17'
$ws.Cells.Item(26, 5).Value = '  +4.49%  '
$ws.Cells.Item(27, 5).Value = '  -0.85%  '
$ws.Cells.Item(28, 4).Value = '''9.53'
$ws.Cells.Item(28, 5).Value = '  -1.49%  '
$ws.Cells.Item(29, 4).Value = '''33.19'
$ws.Cells.Item(29, 5).Value = '  -0.38%  '
$ws.Cells.Item(30, 4).Value = '''7.07'
$ws.Cells.Item(30, 5).Value = '  +1.60%  '
$ws.Cells.Item(31, 4).Value = '''4.30'
$ws.Cells.Item(31, 5).Value = '  +3.60%  '
$ws.Cells.Item(32, 4).Value = '''12.40'
$ws.Cells.Item(32, 5).Value = '  -0.74%  '
$ws.Cells.Item(33, 5).Value = '  -0.12%  '
$ws.Cells.Item(34, 4).Value = '''63.11'
$ws.Cells.Item(34, 5).Value = '  -0.52%  '
$ws.Cells.Item(35, 4).Value = '0.0₃0824'
$ws.Cells.Item(35, 5).Value = '  +6.56%  '
$ws.Cells.Item(36, 4).Value = '3.728.31'
$ws.Cells.Item(36, 5).Value = '  +2.82%  '
$ws.Cells.Item(37, 5).Value = '  -2.49%  '
$ws.Cells.Item(38, 4).Value = '''1.00'
$ws.Cells.Item(39, 5).Value = '  -0.28%  '
$ws.Cells.Item(40, 5).Value = '  -0.69%  '
$ws.Cells.Item(41, 4).Value = '''36.37'
$ws.Cells.Item(41, 5).Value = '  -1.11%  '
$ws.Cells.Item(42, 4).Value = '''483.63'
$ws.Cells.Item(42, 5).Value = '  -3.72%  '
$ws.Cells.Item(43, 5).Value = '  -2.76%  '
$ws.Cells.Item(44, 4).Value = '''0.0454'
$ws.Cells.Item(44, 5).Value = '  -2.05%  '
$ws.Cells.Item(45, 5).Value = '  -1.04%  '
$ws.Cells.Item(46, 5).Value = '  -3.19%  '
$ws.Cells.Item(47, 5).Value = '  -1.29%  '
$ws.Cells.Item(48, 4).Value = '''1.00'
$ws.Cells.Item(48, 5).Value = '  +0.20%  '
$ws.Cells.Item(49, 5).Value = '  -4.33%  '
$ws.Cells.Item(50, 4).Value = '''0.000246'
$ws.Cells.Item(50, 5).Value = '  +2.95%  '
$ws.Cells.Item(51, 5).Value = '  +10.24%  '
